$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Item row (row 2) values to point at the latest auto-generated
# test names, and link the Questionaire form cell (I2) to the new
# supplier-flow test form, clearing its old explicit style so it
# matches the default "Normal" formatting used by the rest of the row.
$ws.Range("A2").Value = "TestForm_11/01/2019-18:02:02"
$ws.Range("B2").Value = "AUTO_TEST_TASK_ON_11/01/2019-10:35:48"
$ws.Range("I2").Value = "TestForm_11/01/2019-16:23:44"
$ws.Range("I2").Style = "Normal"
